# Add the missing "WEAPON_*" translation rows to the localization sheet,
# restyle the table with alternating section banding + bold keys, and
# widen the columns to fit the new (longer) content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New rows of data (rows 16-22). Written column-by-column, row-by-row
#    (A, then B, then C) so the shared-string table is populated in the
#    same order the source workbook uses.
# ---------------------------------------------------------------------
$newRows = @(
    @("WEAPON_CHOICE_MENU",            "Choose your weapon", "Choisis ton arme"),
    @("WEAPON_CHOICE_BUTTON",          "Take",                "Prendre"),
    @("WEAPON_GUN",                    "Gun",                 "Pistolet"),
    @("WEAPON_SWORD",                  "Sword",               "Épée"),
    @("WEAPON_ATTRIBUTE_DAMAGE",       "Damage",              "Dégats"),
    @("WEAPON_ATTRIBUTE_RANGE",        "Range",               "Portée"),
    @("WEAPON_ATTRIBUTE_ATTACK_SPEED", "Attack Speed",        "Vitesse d'attaque")
)

$r = 16
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Column widths (A, B, C all now individually sized).
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 35.33
$ws.Columns("B").ColumnWidth = 34.67
$ws.Columns("C").ColumnWidth = 38.67

# ---------------------------------------------------------------------
# 3. Bold the "key" column for every data row (2-22).
# ---------------------------------------------------------------------
$ws.Range("A2:A22").Font.Bold = $true

# ---------------------------------------------------------------------
# 4. Alternating section-banding fills.
#    "gray" band  -> fgColor #B2B2B2 / bgColor #969696
#    "light" band -> fgColor #EEEEEE / bgColor #FFFFCC
# ---------------------------------------------------------------------
$grayRange  = $ws.Range("A2:C3,A8:C10,A13:C15")
$grayRange.Interior.Color = 11711154
$grayRange.Interior.PatternColor = 9868950

$lightRange = $ws.Range("A4:C7,A11:C12,A16:C22")
$lightRange.Interior.Color = 15658734
$lightRange.Interior.PatternColor = 13434879

# ---------------------------------------------------------------------
# 5. Leave the cursor where the author left it.
# ---------------------------------------------------------------------
$ws.Range("E9").Select()
